$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new methodology run added 7 additional dyes to the dataset (10c, 10d, 12b, 19a, 22a,
# Phenothiazine3, Phenothiazine4) and refreshed the computed properties for all entries.
# Insert 7 blank rows inside the existing table (this also pushes the later rows down and
# automatically extends the sheet dimension from K18 to K25).
$ws.Range("A3:A9").EntireRow.Insert()

# Rewrite every data row (2-25), sorted alphabetically by dye name, with the final values.

# Row 2: 10b
$ws.Cells.Item(2,1).Value = "10b"
$ws.Cells.Item(2,2).Value = -4.597
$ws.Cells.Item(2,3).Value = -3.71
$ws.Cells.Item(2,4).Value = 27.1424
$ws.Cells.Item(2,5).Value = -2857.552262
$ws.Cells.Item(2,6).Value = -1.232
$ws.Cells.Item(2,7).Value = 598.841204
$ws.Cells.Item(2,8).Value = 651.33504
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(2,10).Value = 583
$ws.Cells.Item(2,11).Value = 1.010128

# Row 3: 10c
$ws.Cells.Item(3,1).Value = "10c"
$ws.Cells.Item(3,2).Value = -4.684
$ws.Cells.Item(3,3).Value = -3.682
$ws.Cells.Item(3,4).Value = 20.2171
$ws.Cells.Item(3,5).Value = -2677.852681
$ws.Cells.Item(3,6).Value = -1.105
$ws.Cells.Item(3,7).Value = 763.223679
$ws.Cells.Item(3,8).Value = 840.423948
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(3,10).Value = 558
$ws.Cells.Item(3,11).Value = 1.084045

# Row 4: 10d
$ws.Cells.Item(4,1).Value = "10d"
$ws.Cells.Item(4,2).Value = -4.615
$ws.Cells.Item(4,3).Value = -3.677
$ws.Cells.Item(4,4).Value = 18.4425
$ws.Cells.Item(4,5).Value = -2523.17176
$ws.Cells.Item(4,6).Value = -1.031
$ws.Cells.Item(4,7).Value = 697.440336
$ws.Cells.Item(4,8).Value = 770.669635
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 510
$ws.Cells.Item(4,11).Value = 0.547976

# Row 5: 12a
$ws.Cells.Item(5,1).Value = "12a"
$ws.Cells.Item(5,2).Value = -4.636
$ws.Cells.Item(5,3).Value = -3.698
$ws.Cells.Item(5,4).Value = 27.223
$ws.Cells.Item(5,5).Value = -2857.553234
$ws.Cells.Item(5,6).Value = -1.203
$ws.Cells.Item(5,7).Value = 596.582989
$ws.Cells.Item(5,8).Value = 652.896473
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(5,10).Value = 592
$ws.Cells.Item(5,11).Value = 0.999384

# Row 6: 12b
$ws.Cells.Item(6,1).Value = "12b"
$ws.Cells.Item(6,2).Value = -4.624
$ws.Cells.Item(6,3).Value = -3.677
$ws.Cells.Item(6,4).Value = 14.9789
$ws.Cells.Item(6,5).Value = -3498.950342
$ws.Cells.Item(6,6).Value = -0.944
$ws.Cells.Item(6,7).Value = 582.016851
$ws.Cells.Item(6,8).Value = 641.354295
$ws.Cells.Item(6,9).Value = 0
$ws.Cells.Item(6,10).Value = 528
$ws.Cells.Item(6,11).Value = 0.889745

# Row 7: 19a
$ws.Cells.Item(7,1).Value = "19a"
$ws.Cells.Item(7,2).Value = -4.71
$ws.Cells.Item(7,3).Value = -3.587
$ws.Cells.Item(7,4).Value = 16.8789
$ws.Cells.Item(7,5).Value = -1985.333899
$ws.Cells.Item(7,6).Value = -1.026
$ws.Cells.Item(7,7).Value = 544.193692
$ws.Cells.Item(7,8).Value = 582.522059
$ws.Cells.Item(7,9).Value = 0
$ws.Cells.Item(7,10).Value = 530
$ws.Cells.Item(7,11).Value = 0.71764

# Row 8: 20a
$ws.Cells.Item(8,1).Value = "20a"
$ws.Cells.Item(8,2).Value = -4.651
$ws.Cells.Item(8,3).Value = -3.727
$ws.Cells.Item(8,4).Value = 26.3601
$ws.Cells.Item(8,5).Value = -2216.160004
$ws.Cells.Item(8,6).Value = -1.214
$ws.Cells.Item(8,7).Value = 608.199673
$ws.Cells.Item(8,8).Value = 664.73637
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(8,10).Value = 561
$ws.Cells.Item(8,11).Value = 1.05693

# Row 9: 21a
$ws.Cells.Item(9,1).Value = "21a"
$ws.Cells.Item(9,2).Value = -4.561
$ws.Cells.Item(9,3).Value = -3.715
$ws.Cells.Item(9,4).Value = 23.7772
$ws.Cells.Item(9,5).Value = -2767.699232
$ws.Cells.Item(9,6).Value = -1.31
$ws.Cells.Item(9,7).Value = 678.557941
$ws.Cells.Item(9,8).Value = 748.458472
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,10).Value = 637
$ws.Cells.Item(9,11).Value = 0.705128

# Row 10: 22a
$ws.Cells.Item(10,1).Value = "22a"
$ws.Cells.Item(10,2).Value = -4.558
$ws.Cells.Item(10,3).Value = -3.68
$ws.Cells.Item(10,4).Value = 21.2796
$ws.Cells.Item(10,5).Value = -3319.255172
$ws.Cells.Item(10,6).Value = -1.171
$ws.Cells.Item(10,7).Value = 747.679262
$ws.Cells.Item(10,8).Value = 828.789887
$ws.Cells.Item(10,9).Value = 0
$ws.Cells.Item(10,10).Value = 682
$ws.Cells.Item(10,11).Value = 1.470767

# Row 11: D1
$ws.Cells.Item(11,1).Value = "D1"
$ws.Cells.Item(11,2).Value = -5.195
$ws.Cells.Item(11,3).Value = -2.63
$ws.Cells.Item(11,4).Value = 6.5592
$ws.Cells.Item(11,5).Value = -1393.134805
$ws.Cells.Item(11,6).Value = -1.184
$ws.Cells.Item(11,7).Value = 387.3409
$ws.Cells.Item(11,8).Value = 412.951353
$ws.Cells.Item(11,9).Value = 0
$ws.Cells.Item(11,10).Value = 408
$ws.Cells.Item(11,11).Value = 0.17197

# Row 12: D2
$ws.Cells.Item(12,1).Value = "D2"
$ws.Cells.Item(12,2).Value = -5.163
$ws.Cells.Item(12,3).Value = -2.589
$ws.Cells.Item(12,4).Value = 6.8226
$ws.Cells.Item(12,5).Value = -1278.684862
$ws.Cells.Item(12,6).Value = -1.054
$ws.Cells.Item(12,7).Value = 354.439145
$ws.Cells.Item(12,8).Value = 380.068067
$ws.Cells.Item(12,9).Value = 0
$ws.Cells.Item(12,10).Value = 409
$ws.Cells.Item(12,11).Value = 0.234811

# Row 13: D3
$ws.Cells.Item(13,1).Value = "D3"
$ws.Cells.Item(13,2).Value = -5.271
$ws.Cells.Item(13,3).Value = -2.633
$ws.Cells.Item(13,4).Value = 3.3748
$ws.Cells.Item(13,5).Value = -1444.874262
$ws.Cells.Item(13,6).Value = -1.061
$ws.Cells.Item(13,7).Value = 346.581779
$ws.Cells.Item(13,8).Value = 402.922143
$ws.Cells.Item(13,9).Value = 0
$ws.Cells.Item(13,10).Value = 359
$ws.Cells.Item(13,11).Value = 0.109098

# Row 14: DPAA
$ws.Cells.Item(14,1).Value = "DPAA"
$ws.Cells.Item(14,2).Value = -4.862
$ws.Cells.Item(14,3).Value = -2.735
$ws.Cells.Item(14,4).Value = 5.8789
$ws.Cells.Item(14,5).Value = -3704.69389
$ws.Cells.Item(14,6).Value = -0.685
$ws.Cells.Item(14,7).Value = 389.939641
$ws.Cells.Item(14,8).Value = 426.054576
$ws.Cells.Item(14,9).Value = 0
$ws.Cells.Item(14,10).Value = 474
$ws.Cells.Item(14,11).Value = 0.499187

# Row 15: DPACA
$ws.Cells.Item(15,1).Value = "DPACA"
$ws.Cells.Item(15,2).Value = -5.095
$ws.Cells.Item(15,3).Value = -3.249
$ws.Cells.Item(15,4).Value = 7.7768
$ws.Cells.Item(15,5).Value = -3796.874182
$ws.Cells.Item(15,6).Value = -0.816
$ws.Cells.Item(15,7).Value = 412.305213
$ws.Cells.Item(15,8).Value = 452.543801
$ws.Cells.Item(15,9).Value = 0
$ws.Cells.Item(15,10).Value = 519
$ws.Cells.Item(15,11).Value = 0.533851

# Row 16: Phenothiazine1
$ws.Cells.Item(16,1).Value = "Phenothiazine1"
$ws.Cells.Item(16,2).Value = -4.935
$ws.Cells.Item(16,3).Value = -3.333
$ws.Cells.Item(16,4).Value = 9.0352
$ws.Cells.Item(16,5).Value = -1351.57215
$ws.Cells.Item(16,6).Value = -0.765
$ws.Cells.Item(16,7).Value = 330.452395
$ws.Cells.Item(16,8).Value = 349.95115
$ws.Cells.Item(16,9).Value = 0
$ws.Cells.Item(16,10).Value = 343
$ws.Cells.Item(16,11).Value = 0.288586

# Row 17: Phenothiazine2
$ws.Cells.Item(17,1).Value = "Phenothiazine2"
$ws.Cells.Item(17,2).Value = -4.818
$ws.Cells.Item(17,3).Value = -2.806
$ws.Cells.Item(17,4).Value = 7.2013
$ws.Cells.Item(17,5).Value = -1259.394134
$ws.Cells.Item(17,6).Value = -0.631
$ws.Cells.Item(17,7).Value = 311.25871
$ws.Cells.Item(17,8).Value = 327.396927
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 308
$ws.Cells.Item(17,11).Value = 0.353802

# Row 18: Phenothiazine3
$ws.Cells.Item(18,1).Value = "Phenothiazine3"
$ws.Cells.Item(18,2).Value = -5.34
$ws.Cells.Item(18,3).Value = -3.153
$ws.Cells.Item(18,4).Value = 9.7893
$ws.Cells.Item(18,5).Value = -953.539237
$ws.Cells.Item(18,6).Value = -0.765
$ws.Cells.Item(18,7).Value = 310.115745
$ws.Cells.Item(18,8).Value = 326.248824
$ws.Cells.Item(18,9).Value = 0
$ws.Cells.Item(18,10).Value = 332
$ws.Cells.Item(18,11).Value = 0.479784

# Row 19: Phenothiazine4
$ws.Cells.Item(19,1).Value = "Phenothiazine4"
$ws.Cells.Item(19,2).Value = -4.991
$ws.Cells.Item(19,3).Value = -3.206
$ws.Cells.Item(19,4).Value = 9.4053
$ws.Cells.Item(19,5).Value = -1032.044192
$ws.Cells.Item(19,6).Value = -0.759
$ws.Cells.Item(19,7).Value = 333.827291
$ws.Cells.Item(19,8).Value = 358.719796
$ws.Cells.Item(19,9).Value = 0
$ws.Cells.Item(19,10).Value = 524
$ws.Cells.Item(19,11).Value = 0.534846

# Row 20: Quercitin
$ws.Cells.Item(20,1).Value = "Quercitin"
$ws.Cells.Item(20,2).Value = -5.323
$ws.Cells.Item(20,3).Value = -2.754
$ws.Cells.Item(20,4).Value = 6.6676
$ws.Cells.Item(20,5).Value = -1638.590312
$ws.Cells.Item(20,6).Value = -1.058
$ws.Cells.Item(20,7).Value = 381.893317
$ws.Cells.Item(20,8).Value = 426.986904
$ws.Cells.Item(20,9).Value = 0
$ws.Cells.Item(20,10).Value = 419
$ws.Cells.Item(20,11).Value = 0.193758

# Row 21: Rutin
$ws.Cells.Item(21,1).Value = "Rutin"
$ws.Cells.Item(21,2).Value = -5.159
$ws.Cells.Item(21,3).Value = -2.833
$ws.Cells.Item(21,4).Value = 6.3594
$ws.Cells.Item(21,5).Value = -2710.037461
$ws.Cells.Item(21,6).Value = -1.971
$ws.Cells.Item(21,7).Value = 631.881213
$ws.Cells.Item(21,8).Value = 731.916619
$ws.Cells.Item(21,9).Value = 0
$ws.Cells.Item(21,10).Value = 466
$ws.Cells.Item(21,11).Value = 0.255696

# Row 22: T1
$ws.Cells.Item(22,1).Value = "T1"
$ws.Cells.Item(22,2).Value = -5.418
$ws.Cells.Item(22,3).Value = -3.71
$ws.Cells.Item(22,4).Value = 14.832
$ws.Cells.Item(22,5).Value = -1464.797419
$ws.Cells.Item(22,6).Value = -1.012
$ws.Cells.Item(22,7).Value = 373.33181
$ws.Cells.Item(22,8).Value = 386.599372
$ws.Cells.Item(22,9).Value = 0
$ws.Cells.Item(22,10).Value = 523
$ws.Cells.Item(22,11).Value = 0.952597

# Row 23: T2
$ws.Cells.Item(23,1).Value = "T2"
$ws.Cells.Item(23,2).Value = -5.396
$ws.Cells.Item(23,3).Value = -3.71
$ws.Cells.Item(23,4).Value = 15.8432
$ws.Cells.Item(23,5).Value = -1618.511711
$ws.Cells.Item(23,6).Value = -1.097
$ws.Cells.Item(23,7).Value = 417.523898
$ws.Cells.Item(23,8).Value = 431.785193
$ws.Cells.Item(23,9).Value = 0
$ws.Cells.Item(23,10).Value = 540
$ws.Cells.Item(23,11).Value = 0.794487

# Row 24: T3
$ws.Cells.Item(24,1).Value = "T3"
$ws.Cells.Item(24,2).Value = -5.448
$ws.Cells.Item(24,3).Value = -3.863
$ws.Cells.Item(24,4).Value = 13.091
$ws.Cells.Item(24,5).Value = -1633.444146
$ws.Cells.Item(24,6).Value = -0.952
$ws.Cells.Item(24,7).Value = 393.928755
$ws.Cells.Item(24,8).Value = 413.254515
$ws.Cells.Item(24,9).Value = 0
$ws.Cells.Item(24,10).Value = 560
$ws.Cells.Item(24,11).Value = 1.005572

# Row 25: T4
$ws.Cells.Item(25,1).Value = "T4"
$ws.Cells.Item(25,2).Value = -5.429
$ws.Cells.Item(25,3).Value = -3.866
$ws.Cells.Item(25,4).Value = 13.7502
$ws.Cells.Item(25,5).Value = -1787.160542
$ws.Cells.Item(25,6).Value = -1.021
$ws.Cells.Item(25,7).Value = 441.173419
$ws.Cells.Item(25,8).Value = 457.649836
$ws.Cells.Item(25,9).Value = 0
$ws.Cells.Item(25,10).Value = 583
$ws.Cells.Item(25,11).Value = 0.794271

Write-Host "Updated dataset to A1:K25"
